$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.729.22"
$ws.Range("E2").Value = "  +2.39%  "

$ws.Range("D3").Value = "2.638.02"
$ws.Range("E3").Value = "  +9.52%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.99%  "

$ws.Range("E7").Value = "  +7.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +15.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0846"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.69%  "

$ws.Range("E13").Value = "  +16.49%  "

$ws.Range("D14").Value = "3.037.21"
$ws.Range("E14").Value = "  +10.08%  "

$ws.Range("E15").Value = "  +2.03%  "

$ws.Range("D16").Value = "2.627.99"
$ws.Range("E16").Value = "  +8.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.922"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.05%  "

$ws.Range("D19").Value = "47.084.96"
$ws.Range("E19").Value = "  +3.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.15%  "

$ws.Range("E21").Value = "  +8.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "256.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +15.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +37.28%  "

$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "42.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +16.55%  "

$ws.Range("E35").Value = "  +4.45%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0842"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.36%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.54%  "

$ws.Range("E38").Value = "  +5.03%  "

$ws.Range("E39").Value = "  +6.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0331"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +40.70%  "

$ws.Range("D45").Value = "2.042.69"
$ws.Range("E45").Value = "  +4.67%  "

$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "92.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "114.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.82%  "

$ws.Range("E49").Value = "  +4.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.204"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.11%  "

